$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D/E hold numeric- or percentage-looking text. Excel would silently
# coerce strings like "4.500" or "6.10%" into Numbers (losing the exact
# authored formatting), so force each such cell to Text before assigning it.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "333.08"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "44.26"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "6.10%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.844"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "3.97%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08340"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.09%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "8.819"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "0.97%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.979"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "-2.08%"
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.898"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-2.84%"
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9346"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.40%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1256"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-1.51%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1948"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.40%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09573"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.86%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03935"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "3.28%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1066"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.76%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001304"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.04%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006085"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-3.19%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.506"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.89%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.500"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.56%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.971"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "8.08%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.62%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "6.58%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04413"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "0.01%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001259"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.07%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004407"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "1.57%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.85%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003994"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02797"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "0.88%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05699"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "5.22%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007946"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "3.61%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1425"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.54%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008997"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.15%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002104"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.41%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-10.22%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007272"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "9.36%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.00%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003253"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "1.19%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.12%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.00%"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.00%"
